$wb = $excel.ActiveWorkbook

# --- Sheet2 ("Sheet2 - Numbers"): extend the data with a new column AA ---
# Values 100..129 for rows 1..30 (row r -> 99 + r), matching the "reader/sheet2"
# fixture this test data is being made consistent with.
$ws2 = $wb.Worksheets.Item(2)
for ($r = 1; $r -le 30; $r++) {
    $ws2.Cells.Item($r, 27).Value = 99 + $r
}

# Sheet2 becomes the active sheet/tab, scrolled/selected over the new column.
$ws2.Activate()
$ws2.Range("AA1:AA30").Select() | Out-Null

# --- Sheet4 ("Sheet4 - Dates"): page setup tweak that came along with the resave ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.PageSetup.PaperSize = 9
